$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10 and 11 need their observation data swapped (same field layout,
# different record). Columns D, I, P, S, T, U, V, W, AD, AE, AG, AT, AW,
# AX, AY are identical between the two rows already, so only the columns
# below actually need new values.

$row10 = 10
$row11 = 11

# Force the date/time text columns to stay plain text (they are stored as
# text in the workbook, not real Excel dates) so assigning a date-looking
# string does not get auto-converted into a serial date value.
$ws.Range("Y$($row10):AB$($row11)").NumberFormat = "@"

# Read current (pre-swap) values with Value2 so we get plain scalars
# instead of COM Variant wrappers.
$A10 = $ws.Range("A$row10").Value2
$B10 = $ws.Range("B$row10").Value2
$E10 = $ws.Range("E$row10").Value2
$F10 = $ws.Range("F$row10").Value2
$G10 = $ws.Range("G$row10").Value2
$H10 = $ws.Range("H$row10").Value2
$Q10 = $ws.Range("Q$row10").Value2
$R10 = $ws.Range("R$row10").Value2
$Y10 = $ws.Range("Y$row10").Value2
$Z10 = $ws.Range("Z$row10").Value2
$AA10 = $ws.Range("AA$row10").Value2
$AB10 = $ws.Range("AB$row10").Value2
$AC10 = $ws.Range("AC$row10").Value2

$A11 = $ws.Range("A$row11").Value2
$B11 = $ws.Range("B$row11").Value2
$E11 = $ws.Range("E$row11").Value2
$F11 = $ws.Range("F$row11").Value2
$G11 = $ws.Range("G$row11").Value2
$H11 = $ws.Range("H$row11").Value2
$Q11 = $ws.Range("Q$row11").Value2
$R11 = $ws.Range("R$row11").Value2
$Y11 = $ws.Range("Y$row11").Value2
$Z11 = $ws.Range("Z$row11").Value2
$AA11 = $ws.Range("AA$row11").Value2
$AB11 = $ws.Range("AB$row11").Value2
$AC11 = $ws.Range("AC$row11").Value2

# Write row 10 <- old row 11 values
$ws.Range("A$row10").Value2 = $A11
$ws.Range("B$row10").Value2 = $B11
$ws.Range("E$row10").Value2 = $E11
$ws.Range("F$row10").Value2 = $F11
$ws.Range("G$row10").Value2 = $G11
$ws.Range("H$row10").Value2 = $H11
$ws.Range("Q$row10").Value2 = $Q11
$ws.Range("R$row10").Value2 = $R11
$ws.Range("Y$row10").Value2 = $Y11
$ws.Range("Z$row10").Value2 = $Z11
$ws.Range("AA$row10").Value2 = $AA11
$ws.Range("AB$row10").Value2 = $AB11
$ws.Range("AC$row10").Value2 = $AC11

# Write row 11 <- old row 10 values
$ws.Range("A$row11").Value2 = $A10
$ws.Range("B$row11").Value2 = $B10
$ws.Range("E$row11").Value2 = $E10
$ws.Range("F$row11").Value2 = $F10
$ws.Range("G$row11").Value2 = $G10
$ws.Range("H$row11").Value2 = $H10
$ws.Range("Q$row11").Value2 = $Q10
$ws.Range("R$row11").Value2 = $R10
$ws.Range("Y$row11").Value2 = $Y10
$ws.Range("Z$row11").Value2 = $Z10
$ws.Range("AA$row11").Value2 = $AA10
$ws.Range("AB$row11").Value2 = $AB10
$ws.Range("AC$row11").Value2 = $AC10

$wb.Save()
